# Reorder the "Recorded By" (column G) contributor list on the
# "Session Analysis Results" sheet: wherever "System" appears in the
# comma-separated list but is not already the first entry, move it one
# position earlier (swap it with the entry immediately before it) -
# except when that immediately-preceding entry is "admin@admin.com",
# which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "
    $idx = [Array]::IndexOf($parts, "System")

    if ($idx -le 0) {
        continue
    }

    $predecessor = $parts[$idx - 1]
    if ($predecessor -eq "admin@admin.com") {
        continue
    }

    $tmp = $parts[$idx - 1]
    $parts[$idx - 1] = $parts[$idx]
    $parts[$idx] = $tmp

    $cell.Value = [string]::Join(", ", $parts)
}
